$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Hunk 1: merge the two runs that make up the "SAT Sep 01 ..." line
# into a single run (the two adjacent runs already read as one
# contiguous string, so a find/replace over that exact text collapses
# them into a single run on save).
# -----------------------------------------------------------------
$d.Content.Find.Execute("SAT Sep 01 11:42:32 IST 2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "SAT Sep 01 11:42:32 IST 2018", 2) | Out-Null

# -----------------------------------------------------------------
# Hunk 2: append a brand-new "chick in" entry (17/09/2018 MAMATHA)
# right after the last existing entry ("Amount balance ... - 1144.0"),
# and before the trailing blank paragraphs that close out the document.
# -----------------------------------------------------------------

# Locate the "Amount balance ... - 1144.0" paragraph robustly (it is the
# last populated line of the last existing entry).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Amount balance*1144.0*") {
        $anchorIndex = $i
    }
}

$anchor = $d.Paragraphs($anchorIndex)
$ar = $anchor.Range
$ar.Collapse(0)
$ar.InsertParagraphAfter()

# New paragraph #1: blank, bold (inherits bold from the "Amount balance" mark)
$idx = $anchorIndex + 1

# New paragraph #2: date line
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = 0
$p.Range.InsertBefore("Sat Sep 15")
$d.Paragraphs($idx).Range.InsertAfter(" 13:14:22 IST 2018")

# New paragraph #3: Person Name - KR
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = 0
$p.Range.InsertBefore("Person Name" + "`t" + "`t" + "`t" + "`t" + "- KR")

# New paragraph #4: dashed separator
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = 0
$p.Range.InsertBefore("---------------------------------------------------------------")

# New paragraph #5: Item Name - BEET
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = 0
$p.Range.InsertBefore("Item Name" + "`t" + "`t" + "`t" + "`t" + "- BEET")

# New paragraph #6: Amount Received - 1144 (red)
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = 0
$p.Range.InsertBefore("Amount Received" + "`t" + "`t" + "`t" + "- 1144")
$p.Range.Font.Color = 255

# New paragraph #7: Amount Received mode - CASH AND CLEARD
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = 0
$p.Range.Font.Color = -16777216
$p.Range.InsertBefore("Amount Received mode" + "`t" + "`t" + "- CASH AND CLEARD")

# New paragraph #8: blank, not bold
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = 0

# New paragraph #9: blank, bold
$r = $d.Paragraphs($idx).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs($idx)
$p.Range.Font.Bold = -1

Write-Output "done; paragraphs=$($d.Paragraphs.Count)"
